$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the "Ace hardware" expense amount in I19 ---
$ws.Cells.Item(19, 9).Value = 50.97

# --- Apply the yellow "receipt highlight" fill (same style already used on
#     D16/E17/F18/E27/D28) to the "Other" expense cells I19:I26 ---
$yellow = $ws.Cells.Item(16, 4).Interior.Color
for ($r = 19; $r -le 26; $r++) {
    $ws.Cells.Item($r, 9).Interior.Color = $yellow
}

# --- Note what each "Other" expense line was for, in column K ---
$ws.Cells.Item(19, 11).Value = "Hardware supplies for Puerto Rico field work"
$ws.Cells.Item(22, 11).Value = "electrical supplies for Puerto Rico field work"

# --- Add the two sub-total check formulas in column K ---
$ws.Cells.Item(20, 11).Formula = "=SUM(I19:I21,I23,)"
$ws.Cells.Item(23, 11).Formula = "=SUM(I22,I24,I25,I26)"

# --- Column width touch-up: columns C and F keep their width but lose the
#     autofit ("best fit") flag, and the newly-used D/E/G/H columns pick up
#     an explicit (default) width ---
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 13.333333333333334
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 10

# --- Scroll the sheet up a bit and move the active selection to I22 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 3
$ws.Range("I22").Select()
